$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.141061488777358
$ws.Range("B3").Value = 10.00000009984927
$ws.Range("B4").Value = 65.00000064741765
$ws.Range("B5").Value = 1.256345602973553
$ws.Range("B6").Value = 2.366201640490309
$ws.Range("B7").Value = 0.082776784967591
$ws.Range("B8").Value = 8.087660995730328
$ws.Range("B9").Value = 1.504116959379239
$ws.Range("B10").Value = 1.419258604085065
$ws.Range("B11").Value = -0.08485835529417374
$ws.Range("B12").Value = 0.03999999054947855
$ws.Range("B13").Value = 0.3999999900529826
$ws.Range("B14").Value = 0.01999999000524019
$ws.Range("B15").Value = 49463.23004825977
$ws.Range("B16").Value = -0.1015359300736017
$ws.Range("B17").Value = 0.8056713411417376
$ws.Range("B18").Value = 0.8114614932495552
$ws.Range("B19").Value = 2.053750222914737
$ws.Range("B20").Value = 0.1405388941806902
$ws.Range("B22").Value = 0.1140419009315041
$ws.Range("B23").Value = 0.8056713411417376
$ws.Range("B24").Value = 0.1784068470628501
$ws.Range("B25").Value = 0.04460171176571253
$ws.Range("B27").Value = 0.007957250769726797
$ws.Range("B28").Value = 0.8056713411417376
$ws.Range("B29").Value = 0.06778242641157563
$ws.Range("B30").Value = 0.03389121320578781
$ws.Range("B32").Value = 0.002297228665120333
$ws.Range("B36").Value = [double]"1.873081839836726e-16"
$ws.Range("B37").Value = [double]"1.051024499184463e-14"
$ws.Range("B38").Value = 10.00580694703766
$ws.Range("B39").Value = 0.2288781246149695
$ws.Range("B40").Value = 1.363862811138331
$ws.Range("B41").Value = 4.956629328400594
$ws.Range("B42").Value = -40.62063333904845
$ws.Range("B43").Value = 0.4866281482450332
$ws.Range("B44").Value = 1.065901365063418
$ws.Range("B45").Value = 0.007249436612822757
$ws.Range("B46").Value = 1.057650766723981
$ws.Range("B47").Value = 0.001001161726613871
$ws.Range("B48").Value = [double]"3.425414303793734e-10"
$ws.Range("B49").Value = 0.03284266895617405
$ws.Range("B50").Value = [double]"1.657703153232167e-20"
$ws.Range("B51").Value = 0.003182731310834771
$ws.Range("B52").Value = 0.001945924981551996
$ws.Range("B53").Value = 0.001961788972003531
$ws.Range("B54").Value = 0.003903935295313687
$ws.Range("B55").Value = [double]"-1.02986929694753e-21"
$ws.Range("B56").Value = [double]"-4.487751216942977e-06"
$ws.Range("B57").Value = [double]"-2.741783423357421e-21"
$ws.Range("B58").Value = 0.03289764425245234
$ws.Range("B59").Value = 0.0006418656542215293
$ws.Range("B60").Value = 0.0001408059179581332
$ws.Range("B65").Value = 0.008164865440938052
$ws.Range("B67").Value = [double]"-9.999282099889604e-09"
$ws.Range("B68").Value = 0.006131613701303038
$ws.Range("B69").Value = -0.03513472354517255
$ws.Range("B72").Value = -0.03513472354517255
$ws.Range("B73").Value = 0.4054701244248722
$ws.Range("B74").Value = 0.0145401841658935
$ws.Range("B75").Value = 0.7610696293760251
$ws.Range("B78").Value = 0.7610696293760251
$ws.Range("B79").Value = 0.08920342353142506
$ws.Range("B81").Value = 0.7717801279359497
$ws.Range("B84").Value = 0.7717801279359497
$ws.Range("B86").Value = 0.06778242641157563
$ws.Range("B87").Value = -0.09968772864805175
$ws.Range("B88").Value = [double]"-5.16967716399672e-18"
$ws.Range("B89").Value = 1.933514829472097
$ws.Range("B90").Value = [double]"1.393605914562689e-19"
$ws.Range("B91").Value = [double]"2.369068483054626e-14"
$ws.Range("B92").Value = [double]"-3.743027284758309e-18"
$ws.Range("B93").Value = 0.09968772864805175
$ws.Range("B94").Value = -0.2375780819255464
$ws.Range("B95").Value = [double]"-1.393605914562689e-19"
$ws.Range("B96").Value = [double]"5.127264824914205e-19"
$ws.Range("B97").Value = 1.921450986629662
$ws.Range("B98").Value = [double]"-5.16967716399672e-18"
$ws.Range("B99").Value = 0.2375780819255464
$ws.Range("B100").Value = [double]"-1.393605914562689e-19"
$ws.Range("B101").Value = [double]"2.369068483054626e-14"
$ws.Range("B102").Value = [double]"3.743027284758309e-18"
$ws.Range("B103").Value = 1.04076690245188
$ws.Range("B104").Value = [double]"-2.80019054666946e-18"
$ws.Range("B105").Value = 0.1286857724379556
$ws.Range("B106").Value = [double]"-9.302426065963594e-20"
$ws.Range("B107").Value = [double]"9.13072462700741e-14"
$ws.Range("B108").Value = [double]"2.498499340129082e-18"
$ws.Range("B109").Value = 1.85433148830393
$ws.Range("B110").Value = 0.0687012905396126
$ws.Range("B113").Value = 0.05210876440430913
$ws.Range("B114").Value = -0
$ws.Range("B115").Value = 0.810940416261092
$ws.Range("B116").Value = 0.799863895922665
$ws.Range("B117").Value = 0.06518784432642251
$ws.Range("B118").Value = 0.02827084386690126
$ws.Range("B120").Value = -0
$ws.Range("B121").Value = -0.05133653947947021
$ws.Range("B122").Value = -0
$ws.Range("B123").Value = 0.1784068470628501
$ws.Range("B125").Value = -0
$ws.Range("B126").Value = 0.001118499802978656
$ws.Range("B127").Value = [double]"-5.16967716399672e-18"
$ws.Range("B128").Value = [double]"-1.393963855822461e-19"
$ws.Range("B129").Value = -0.0001109254283923663
$ws.Range("B130").Value = [double]"3.743045881230322e-18"
$ws.Range("B131").Value = 0.06778242641157563
$ws.Range("B133").Value = 0.001931653999308799
$ws.Range("B134").Value = 0.001479730226180829
$ws.Range("B136").Value = [double]"3.579412597724902e-23"
$ws.Range("B137").Value = -0.0006612994964228587
$ws.Range("B138").Value = [double]"-1.859647201169946e-23"
$ws.Range("B141").Value = 0.09957036443567334
$ws.Range("B142").Value = 0.1380077174898731
$ws.Range("B143").Value = 4.021942757509349
$ws.Range("B144").Value = 0.7196226212683181
$ws.Range("B145").Value = [double]"1.044407914869513e-22"
$ws.Range("B146").Value = [double]"7.352052803835053e-20"
$ws.Range("B147").Value = -0.1608776722909176
$ws.Range("B148").Value = [double]"1.283184935960032e-20"
$ws.Range("B149").Value = 0.03846422339523582
$ws.Range("B150").Value = [double]"-5.438672171990259e-06"
$ws.Range("B151").Value = [double]"-6.659143239353572e-07"
$ws.Range("B152").Value = 0.1396629135786665
$ws.Range("B153").Value = -0.02499999013591874
$ws.Range("B154").Value = [double]"6.800997353113348e-07"
$ws.Range("B155").Value = 0.03217912957886197
$ws.Range("B156").Value = -0.1541226740289295
$ws.Range("B157").Value = -0.008244218856257035
$ws.Range("B158").Value = -0.001687347275042273
$ws.Range("B159").Value = -0.07811761688372774
$ws.Range("B160").Value = -0.6462573293105381
$ws.Range("B161").Value = -0.00121457648932765
$ws.Range("B162").Value = -0.2729271052266445
$ws.Range("B163").Value = 0.7477096698742525
$ws.Range("B164").Value = 1.261830316836743
$ws.Range("B165").Value = [double]"-8.218579486610794e-18"
$ws.Range("B166").Value = [double]"5.335777174684998e-19"
$ws.Range("B167").Value = -16.67366636131237
$ws.Range("B168").Value = [double]"7.47333886625152e-18"
$ws.Range("B169").Value = -0.0001276918790704684
$ws.Range("B170").Value = [double]"-3.214774066395698e-05"
$ws.Range("B171").Value = 0.1128903331558799
$ws.Range("B172").Value = 0.2225041885802361
$ws.Range("B173").Value = -0.0002291904815012994
$ws.Range("B174").Value = -0.07231340753291615
